# sprint 108 MT cases
# Row heights on the manual-testcases sheet grew (wrapped text got taller)
# for the rows touched in this sprint, and the view should no longer be
# scrolled down to row 6 - it should show from the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(9).RowHeight  = 48.75
$ws.Rows(10).RowHeight = 50.25
$ws.Rows(11).RowHeight = 74.25
$ws.Rows(12).RowHeight = 63.75
$ws.Rows(13).RowHeight = 33
$ws.Rows(14).RowHeight = 43.5
$ws.Rows(15).RowHeight = 45
$ws.Rows(18).RowHeight = 36.75

# Reset the sheet view so it is no longer scrolled to A6 (topLeftCell removed).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
